$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "A working RFID tag" -> "A working RFID tags", split across two runs
#    ("A working RFID " and "tags") with identical run formatting.
# ---------------------------------------------------------------------------
$rfid = $d.Content
[void]$rfid.Find.Execute("A working RFID tag", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)

$wordEnd = $rfid.End
$tagOnly = $d.Range($wordEnd - 3, $wordEnd)   # the "tag" substring
$tagOnly.Text = "tags"

# Force the run boundary to stick (otherwise two adjacent, identically
# formatted runs get coalesced back into one on save) by touching a
# no-op formatting property on just the newly written word.
$splitPoint = $wordEnd - 3
$newWordRange = $d.Range($splitPoint, $splitPoint + 4)
$newWordRange.Font.Bold = 1
$newWordRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) Re-anchor the "_GoBack" bookmark to just after the Postman paragraph
#    text (this is also where the doc's only _Hlk124364050 bookmark IDs get
#    renumbered 0 -> 1, and the old _GoBack bookmark around the UML image
#    gets removed automatically, matching the rest of the diff).
# ---------------------------------------------------------------------------
$postman = $d.Content
[void]$postman.Find.Execute("troubleshooting web API", $true, $false, $false, $false, `
                       $false, $true, 1, $false, "", 0)

$para = $postman.Paragraphs(1)
$paraEnd = $para.Range.End - 1   # position right before the paragraph mark

# Temporarily insert a one-character placeholder so we have a non-empty
# Range to anchor the bookmark to, then delete it, leaving a collapsed
# bookmark exactly at that text position (mirrors the real Word _GoBack
# bookmark, which marks the end point of the most recent edit).
$beforeEnd = $d.Range($paraEnd - 1, $paraEnd)
$beforeEnd.InsertAfter("X")

$placeholder = $d.Range($paraEnd, $paraEnd + 1)
$d.Bookmarks.Add("_GoBack", $placeholder)

$placeholder2 = $d.Range($paraEnd, $paraEnd + 1)
$placeholder2.Delete()
